# Auto-generated edit script for horarios-141-completo.xlsx
# Adds new scraped schedule rows to the three sheets (LP1912, LP1912-215, 6203-6173)
# and refreshes the 'Última actualización' / 'Total filas' header cells.

$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# ---- Sheet "LP1912" ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range('A2').Value = 'Última actualización: 30/12/2025 10:54:50'
$ws1.Range('A3').Value = 'Total filas: 152'

Set-Cell $ws1 132 2 '10:54:39'; Set-Cell $ws1 132 3 '10:56'; Set-Cell $ws1 132 4 '16_SANTA ANA'; Set-Cell $ws1 132 5 2; Set-Cell $ws1 132 6 'LP1912'; Set-Cell $ws1 132 7 '30/12/2025'
Set-Cell $ws1 133 2 '10:54:39'; Set-Cell $ws1 133 3 '10:56'; Set-Cell $ws1 133 4 '27_EL RETIRO'; Set-Cell $ws1 133 5 2; Set-Cell $ws1 133 6 'LP1912'; Set-Cell $ws1 133 7 '30/12/2025'
Set-Cell $ws1 134 2 '10:54:39'; Set-Cell $ws1 134 3 '11:01'; Set-Cell $ws1 134 4 '215C_EL PATO'; Set-Cell $ws1 134 5 7; Set-Cell $ws1 134 6 'LP1912'; Set-Cell $ws1 134 7 '30/12/2025'
Set-Cell $ws1 135 2 '10:54:39'; Set-Cell $ws1 135 3 '11:04'; Set-Cell $ws1 135 4 '23_HERNANDEZ'; Set-Cell $ws1 135 5 10; Set-Cell $ws1 135 6 'LP1912'; Set-Cell $ws1 135 7 '30/12/2025'
Set-Cell $ws1 136 2 '10:54:39'; Set-Cell $ws1 136 3 '11:06'; Set-Cell $ws1 136 4 '16_P MOR-167 Y 521'; Set-Cell $ws1 136 5 12; Set-Cell $ws1 136 6 'LP1912'; Set-Cell $ws1 136 7 '30/12/2025'
Set-Cell $ws1 137 2 '10:54:39'; Set-Cell $ws1 137 3 '11:11'; Set-Cell $ws1 137 4 '10_OLMOS'; Set-Cell $ws1 137 5 17; Set-Cell $ws1 137 6 'LP1912'; Set-Cell $ws1 137 7 '30/12/2025'
Set-Cell $ws1 138 2 '10:54:39'; Set-Cell $ws1 138 3 '11:21'; Set-Cell $ws1 138 4 '26_HERNANDEZ'; Set-Cell $ws1 138 5 27; Set-Cell $ws1 138 6 'LP1912'; Set-Cell $ws1 138 7 '30/12/2025'
Set-Cell $ws1 139 2 '10:54:39'; Set-Cell $ws1 139 3 '11:22'; Set-Cell $ws1 139 4 '10_OLMOS'; Set-Cell $ws1 139 5 28; Set-Cell $ws1 139 6 'LP1912'; Set-Cell $ws1 139 7 '30/12/2025'
Set-Cell $ws1 140 2 '10:54:39'; Set-Cell $ws1 140 3 '11:26'; Set-Cell $ws1 140 4 '16_SANTA ANA'; Set-Cell $ws1 140 5 32; Set-Cell $ws1 140 6 'LP1912'; Set-Cell $ws1 140 7 '30/12/2025'
Set-Cell $ws1 141 2 '10:54:39'; Set-Cell $ws1 141 3 '11:34'; Set-Cell $ws1 141 4 '23_HERNANDEZ'; Set-Cell $ws1 141 5 40; Set-Cell $ws1 141 6 'LP1912'; Set-Cell $ws1 141 7 '30/12/2025'
Set-Cell $ws1 142 2 '10:54:39'; Set-Cell $ws1 142 3 '11:36'; Set-Cell $ws1 142 4 '16_SANTA ANA'; Set-Cell $ws1 142 5 42; Set-Cell $ws1 142 6 'LP1912'; Set-Cell $ws1 142 7 '30/12/2025'
Set-Cell $ws1 143 2 '10:54:39'; Set-Cell $ws1 143 3 '11:42'; Set-Cell $ws1 143 4 '17_ROMERO'; Set-Cell $ws1 143 5 48; Set-Cell $ws1 143 6 'LP1912'; Set-Cell $ws1 143 7 '30/12/2025'
Set-Cell $ws1 144 2 '10:54:39'; Set-Cell $ws1 144 3 '11:43'; Set-Cell $ws1 144 4 '10_OLMOS'; Set-Cell $ws1 144 5 49; Set-Cell $ws1 144 6 'LP1912'; Set-Cell $ws1 144 7 '30/12/2025'
Set-Cell $ws1 145 2 '10:54:39'; Set-Cell $ws1 145 3 '11:52'; Set-Cell $ws1 145 4 '15_ABASTO'; Set-Cell $ws1 145 5 58; Set-Cell $ws1 145 6 'LP1912'; Set-Cell $ws1 145 7 '30/12/2025'
Set-Cell $ws1 146 2 '10:54:39'; Set-Cell $ws1 146 3 '12:02'; Set-Cell $ws1 146 4 '84_COLONIA URQUIZA-ESC 49'; Set-Cell $ws1 146 5 68; Set-Cell $ws1 146 6 'LP1912'; Set-Cell $ws1 146 7 '30/12/2025'
Set-Cell $ws1 147 2 '10:54:39'; Set-Cell $ws1 147 3 '12:06'; Set-Cell $ws1 147 4 '16_P MOR-SANTA ANA'; Set-Cell $ws1 147 5 72; Set-Cell $ws1 147 6 'LP1912'; Set-Cell $ws1 147 7 '30/12/2025'
Set-Cell $ws1 148 2 '10:54:39'; Set-Cell $ws1 148 3 '12:08'; Set-Cell $ws1 148 4 '23_HERNANDEZ'; Set-Cell $ws1 148 5 74; Set-Cell $ws1 148 6 'LP1912'; Set-Cell $ws1 148 7 '30/12/2025'
Set-Cell $ws1 149 2 '10:54:39'; Set-Cell $ws1 149 3 '12:21'; Set-Cell $ws1 149 4 '14_ABASTO'; Set-Cell $ws1 149 5 87; Set-Cell $ws1 149 6 'LP1912'; Set-Cell $ws1 149 7 '30/12/2025'
Set-Cell $ws1 150 2 '10:54:39'; Set-Cell $ws1 150 3 '12:21'; Set-Cell $ws1 150 4 '26_HERNANDEZ'; Set-Cell $ws1 150 5 87; Set-Cell $ws1 150 6 'LP1912'; Set-Cell $ws1 150 7 '30/12/2025'
Set-Cell $ws1 151 2 '10:54:39'; Set-Cell $ws1 151 3 '12:23'; Set-Cell $ws1 151 4 '17_ROMERO'; Set-Cell $ws1 151 5 89; Set-Cell $ws1 151 6 'LP1912'; Set-Cell $ws1 151 7 '30/12/2025'
Set-Cell $ws1 152 2 '10:54:39'; Set-Cell $ws1 152 3 '12:24'; Set-Cell $ws1 152 4 '15_ABASTO'; Set-Cell $ws1 152 5 90; Set-Cell $ws1 152 6 'LP1912'; Set-Cell $ws1 152 7 '30/12/2025'
Set-Cell $ws1 153 2 '10:54:39'; Set-Cell $ws1 153 3 '12:27'; Set-Cell $ws1 153 4 '15_ABASTO'; Set-Cell $ws1 153 5 93; Set-Cell $ws1 153 6 'LP1912'; Set-Cell $ws1 153 7 '30/12/2025'

# ---- Sheet "LP1912-215" ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range('A2').Value = 'Última actualización: 30/12/2025 10:54:50'
$ws2.Range('A3').Value = 'Total filas: 17'

Set-Cell $ws2 18 2 '30/12/2025'; Set-Cell $ws2 18 3 '10:54:39'; Set-Cell $ws2 18 4 '11:01'; Set-Cell $ws2 18 5 '215C_EL PATO'; Set-Cell $ws2 18 6 7; Set-Cell $ws2 18 7 'LP1912'

# ---- Sheet "6203-6173" ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range('A2').Value = 'Última actualización: 30/12/2025 10:54:50'
$ws3.Range('A3').Value = 'Total filas: 21'

Set-Cell $ws3 21 2 '30/12/2025'; Set-Cell $ws3 21 3 '10:54:45'; Set-Cell $ws3 21 4 '11:13'; Set-Cell $ws3 21 5 '215C_LA PLATA'; Set-Cell $ws3 21 6 19; Set-Cell $ws3 21 7 'L6203'
Set-Cell $ws3 22 2 '30/12/2025'; Set-Cell $ws3 22 3 '10:54:50'; Set-Cell $ws3 22 4 '12:04'; Set-Cell $ws3 22 5 '215A_LA PLATA'; Set-Cell $ws3 22 6 70; Set-Cell $ws3 22 7 'L6173'

